$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new header labels for columns I (I0) and J (IF)
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Copy the header style/formatting from H1 (existing header) onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill data rows 2-24:
#   column I (I0) = constant 1
#   column J (IF) = same value as column H (IP) on that row
for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $ws.Cells.Item($r, 8).Value2
}
